$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 108; existing rows 108:129 shift down to 109:130.
$ws.Rows("108:108").Insert()

# Populate the newly inserted row 108 with the new weekly record.
$ws.Range("A108").Value = 8
$ws.Range("B108").Value = "Terminal La Palmera de La Serena"
$ws.Range("C108").Value = "Coquimbo"
$ws.Range("D108").Value = 44637
$ws.Range("E108").Value = 4
$ws.Range("F108").Value = 100112040
$ws.Range("G108").Value = "Cilantro"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 2400
$ws.Range("K108").Value = 2500
$ws.Range("L108").Value = 3000
$ws.Range("M108").Value = 2750
$ws.Range("N108").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O108").Value = "Provincia del Elquí"
$ws.Range("P108").Value = 1833
$ws.Range("Q108").Value = 1.5
$ws.Range("R108").Value = "Hortaliza"
